# Update "想去人数" (column F) counts across sheets, matching the
# gh-pages data refresh (commit 456a3b4). Column F holds a plain
# numeric count per row, so we just overwrite the cell values on the
# relevant worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1129
$ws1.Range("F9").Value = 153
$ws1.Range("F12").Value = 601
$ws1.Range("F13").Value = 1746
$ws1.Range("F14").Value = 1794
$ws1.Range("F15").Value = 842
$ws1.Range("F16").Value = 269
$ws1.Range("F23").Value = 441
$ws1.Range("F24").Value = 86
$ws1.Range("F25").Value = 4731
$ws1.Range("F27").Value = 568
$ws1.Range("F30").Value = 117

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 51
$ws2.Range("F8").Value = 25

# Sheet "全部类型" (all types, aggregated view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 51
$ws4.Range("F11").Value = 25
$ws4.Range("F12").Value = 1129
$ws4.Range("F17").Value = 153
$ws4.Range("F18").Value = 153
$ws4.Range("F22").Value = 601
$ws4.Range("F23").Value = 1746
$ws4.Range("F24").Value = 1794
$ws4.Range("F25").Value = 842
$ws4.Range("F26").Value = 269
$ws4.Range("F35").Value = 441
$ws4.Range("F36").Value = 86
$ws4.Range("F37").Value = 4731
$ws4.Range("F39").Value = 568
$ws4.Range("F44").Value = 117
